$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.31848966666667
$ws.Range("H2").Value = 63.955469
$ws.Range("I2").Value = 0.5519683995553906
$ws.Range("J2").Value = 0.5519683995553906
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7715313333333333
$ws.Range("N2").Value = 2.314594
$ws.Range("O2").Value = 0.05172308417778351
$ws.Range("P2").Value = 0.05172308417778351
$ws.Range("Q2").Value = 16.44788275717622
$ws.Range("R2").Value = 148.030944814586
$ws.Range("S2").Value = 0.02854950799367991
$ws.Range("T2").Value = 0.02854950799367991
$ws.Range("G3").Value = 21.31848966666667
$ws.Range("H3").Value = 63.955469
$ws.Range("I3").Value = 0.5519683995553906
$ws.Range("J3").Value = 0.5519683995553906
$ws.Range("M3").Value = 3.713472666666666
$ws.Range("O3").Value = 0.2489493958723191
$ws.Range("P3").Value = 0.2489493958723191
$ws.Range("Q3").Value = 79.16562867178244
$ws.Range("R3").Value = 712.490658046042
$ws.Range("S3").Value = 0.1374121996099253
$ws.Range("T3").Value = 0.1374121996099253
$ws.Range("G4").Value = 21.31848966666667
$ws.Range("H4").Value = 63.955469
$ws.Range("I4").Value = 0.5519683995553906
$ws.Range("J4").Value = 0.5519683995553906
$ws.Range("M4").Value = 1.701035666666667
$ws.Range("N4").Value = 5.103107
$ws.Range("O4").Value = 0.1140366012048922
$ws.Range("P4").Value = 0.1140366012048922
$ws.Range("Q4").Value = 36.26351128246478
$ws.Range("R4").Value = 326.371601542183
$ws.Range("S4").Value = 0.06294460025780067
$ws.Range("T4").Value = 0.06294460025780067
$ws.Range("G5").Value = 21.31848966666667
$ws.Range("H5").Value = 63.955469
$ws.Range("I5").Value = 0.5519683995553906
$ws.Range("J5").Value = 0.5519683995553906
$ws.Range("M5").Value = 8.730536666666668
$ws.Range("N5").Value = 26.19161
$ws.Range("O5").Value = 0.5852909187450052
$ws.Range("P5").Value = 0.5852909187450052
$ws.Range("Q5").Value = 186.1218557127878
$ws.Range("R5").Value = 1675.09670141509
$ws.Range("S5").Value = 0.3230620916939846
$ws.Range("T5").Value = 0.3230620916939846
$ws.Range("G6").Value = 2.417539666666667
$ws.Range("H6").Value = 7.252619
$ws.Range("I6").Value = 0.06259381042167039
$ws.Range("J6").Value = 0.06259381042167038
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7715313333333333
$ws.Range("N6").Value = 2.314594
$ws.Range("O6").Value = 0.05172308417778351
$ws.Range("P6").Value = 0.05172308417778351
$ws.Range("Q6").Value = 1.865207602409556
$ws.Range("R6").Value = 16.786868421686
$ws.Range("S6").Value = 0.00323754492544828
$ws.Range("T6").Value = 0.00323754492544828
$ws.Range("G7").Value = 2.417539666666667
$ws.Range("H7").Value = 7.252619
$ws.Range("I7").Value = 0.06259381042167039
$ws.Range("J7").Value = 0.06259381042167038
$ws.Range("M7").Value = 3.713472666666666
$ws.Range("O7").Value = 0.2489493958723191
$ws.Range("P7").Value = 0.2489493958723191
$ws.Range("S7").Value = 0.01558269128982131
$ws.Range("T7").Value = 0.01558269128982131
$ws.Range("G8").Value = 2.417539666666667
$ws.Range("H8").Value = 7.252619
$ws.Range("I8").Value = 0.06259381042167039
$ws.Range("J8").Value = 0.06259381042167038
$ws.Range("M8").Value = 1.701035666666667
$ws.Range("N8").Value = 5.103107
$ws.Range("O8").Value = 0.1140366012048922
$ws.Range("P8").Value = 0.1140366012048922
$ws.Range("Q8").Value = 4.112321198581444
$ws.Range("R8").Value = 37.010890787233
$ws.Range("S8").Value = 0.007137985396950652
$ws.Range("T8").Value = 0.007137985396950651
$ws.Range("G9").Value = 2.417539666666667
$ws.Range("H9").Value = 7.252619
$ws.Range("I9").Value = 0.06259381042167039
$ws.Range("J9").Value = 0.06259381042167038
$ws.Range("M9").Value = 8.730536666666668
$ws.Range("N9").Value = 26.19161
$ws.Range("O9").Value = 0.5852909187450052
$ws.Range("P9").Value = 0.5852909187450052
$ws.Range("Q9").Value = 21.10641870295445
$ws.Range("R9").Value = 189.95776832659
$ws.Range("S9").Value = 0.03663558880945014
$ws.Range("T9").Value = 0.03663558880945014
$ws.Range("G10").Value = 0.7420966666666667
$ws.Range("H10").Value = 2.22629
$ws.Range("I10").Value = 0.01921402105965591
$ws.Range("J10").Value = 0.01921402105965591
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7715313333333333
$ws.Range("N10").Value = 2.314594
$ws.Range("O10").Value = 0.05172308417778351
$ws.Range("P10").Value = 0.05172308417778351
$ws.Range("Q10").Value = 0.5725508306955556
$ws.Range("R10").Value = 5.15295747626
$ws.Range("S10").Value = 0.0009938084286622874
$ws.Range("T10").Value = 0.0009938084286622876
$ws.Range("G11").Value = 0.7420966666666667
$ws.Range("H11").Value = 2.22629
$ws.Range("I11").Value = 0.01921402105965591
$ws.Range("J11").Value = 0.01921402105965591
$ws.Range("M11").Value = 3.713472666666666
$ws.Range("O11").Value = 0.2489493958723191
$ws.Range("P11").Value = 0.2489493958723191
$ws.Range("Q11").Value = 2.755755687691111
$ws.Range("R11").Value = 24.80180118922
$ws.Range("S11").Value = 0.004783318935079354
$ws.Range("T11").Value = 0.004783318935079354
$ws.Range("G12").Value = 0.7420966666666667
$ws.Range("H12").Value = 2.22629
$ws.Range("I12").Value = 0.01921402105965591
$ws.Range("J12").Value = 0.01921402105965591
$ws.Range("M12").Value = 1.701035666666667
$ws.Range("N12").Value = 5.103107
$ws.Range("O12").Value = 0.1140366012048922
$ws.Range("P12").Value = 0.1140366012048922
$ws.Range("Q12").Value = 1.262332898114445
$ws.Range("R12").Value = 11.36099608303
$ws.Range("S12").Value = 0.002191101657122381
$ws.Range("T12").Value = 0.002191101657122381
$ws.Range("G13").Value = 0.7420966666666667
$ws.Range("H13").Value = 2.22629
$ws.Range("I13").Value = 0.01921402105965591
$ws.Range("J13").Value = 0.01921402105965591
$ws.Range("M13").Value = 8.730536666666668
$ws.Range("N13").Value = 26.19161
$ws.Range("O13").Value = 0.5852909187450052
$ws.Range("P13").Value = 0.5852909187450052
$ws.Range("Q13").Value = 6.478902158544446
$ws.Range("R13").Value = 58.3101194269
$ws.Range("S13").Value = 0.01124579203879188
$ws.Range("T13").Value = 0.01124579203879188
$ws.Range("G14").Value = 14.14453733333333
$ws.Range("H14").Value = 42.433612
$ws.Range("I14").Value = 0.3662237689632831
$ws.Range("J14").Value = 0.3662237689632831
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.7715313333333333
$ws.Range("N14").Value = 2.314594
$ws.Range("O14").Value = 0.05172308417778351
$ws.Range("P14").Value = 0.05172308417778351
$ws.Range("Q14").Value = 10.91295374816978
$ws.Range("R14").Value = 98.216583733528
$ws.Range("S14").Value = 0.01894222282999303
$ws.Range("T14").Value = 0.01894222282999303
$ws.Range("G15").Value = 14.14453733333333
$ws.Range("H15").Value = 42.433612
$ws.Range("I15").Value = 0.3662237689632831
$ws.Range("J15").Value = 0.3662237689632831
$ws.Range("M15").Value = 3.713472666666666
$ws.Range("O15").Value = 0.2489493958723191
$ws.Range("P15").Value = 0.2489493958723191
$ws.Range("Q15").Value = 52.52535276997954
$ws.Range("R15").Value = 472.7281749298159
$ws.Range("S15").Value = 0.09117118603749308
$ws.Range("T15").Value = 0.09117118603749308
$ws.Range("G16").Value = 14.14453733333333
$ws.Range("H16").Value = 42.433612
$ws.Range("I16").Value = 0.3662237689632831
$ws.Range("J16").Value = 0.3662237689632831
$ws.Range("M16").Value = 1.701035666666667
$ws.Range("N16").Value = 5.103107
$ws.Range("O16").Value = 0.1140366012048922
$ws.Range("P16").Value = 0.1140366012048922
$ws.Range("Q16").Value = 24.06036249249822
$ws.Range("R16").Value = 216.543262432484
$ws.Range("S16").Value = 0.04176291389301849
$ws.Range("T16").Value = 0.04176291389301849
$ws.Range("G17").Value = 14.14453733333333
$ws.Range("H17").Value = 42.433612
$ws.Range("I17").Value = 0.3662237689632831
$ws.Range("J17").Value = 0.3662237689632831
$ws.Range("M17").Value = 8.730536666666668
$ws.Range("N17").Value = 26.19161
$ws.Range("O17").Value = 0.5852909187450052
$ws.Range("P17").Value = 0.5852909187450052
$ws.Range("Q17").Value = 123.4894018217022
$ws.Range("R17").Value = 1111.40461639532
$ws.Range("S17").Value = 0.2143474462027785
$ws.Range("T17").Value = 0.2143474462027785
